# Se verifica que el archivo a procesar exista
# This script:
#  1. Inserts a new worksheet "Resumen total" right after "Punta" (before "Resumen Base")
#     and populates it as a combined summary (built from "Resumen punta" as a template,
#     since it already shares most values with the new sheet).
#  2. Inserts a new row (row 3) into "Resumen Base", "Resumen Intermedia" and "Resumen punta",
#     shifting the existing rows down, and fills the new row with its own
#     Demanda Maxima/Minima entry (index "1", and a consumo total value).

$wb = $excel.ActiveWorkbook

function Set-NewRowThree {
    param($SheetName, $D3Value)

    $ws = $wb.Worksheets.Item($SheetName)

    # Shift existing rows 3..N down by inserting a blank row at position 3.
    $ws.Rows.Item(3).Insert()

    # Recreate the blank placeholder cells for columns B and C (matching the
    # empty-but-present cells used elsewhere in row 2, e.g. B2/C2).
    $ws.Range("B2:C2").Copy($ws.Range("B3:C3"))

    # A3 uses the same bold/centered/bordered style as the rest of column A.
    $a3 = $ws.Range("A3")
    $a3.Font.Bold = $true
    $a3.HorizontalAlignment = -4108
    $a3.VerticalAlignment = -4160
    $a3.Borders.LineStyle = 1
    $a3.Value = 1

    $ws.Range("D3").Value = $D3Value
}

# ---------------------------------------------------------------------------
# 1) Create the new "Resumen total" sheet, positioned right after "Punta".
# ---------------------------------------------------------------------------
$puntaSheet = $wb.Worksheets.Item("Punta")
$totalSheet = $wb.Worksheets.Add($null, $puntaSheet)
$totalSheet.Name = "Resumen total"

# Use the existing "Resumen punta" sheet (closest match) as a formatting/value
# template, then patch the handful of cells that actually differ.
$puntaResumen = $wb.Worksheets.Item("Resumen punta")
$totalSheet2 = $wb.Worksheets.Item("Resumen total")
$puntaResumen.Range("A1:D8").Copy($totalSheet2.Range("A1:D8"))
# The source sheet has no value in A1 (only B1:D1 hold the header labels);
# the range copy still materializes an empty A1 cell, so drop it again.
$totalSheet2.Range("A1").ClearContents()

# ---------------------------------------------------------------------------
# 2) Insert the new "Demanda..." row (row 3) in the brand new sheet and set
#    its own consumo-total value.
# ---------------------------------------------------------------------------
Set-NewRowThree "Resumen total" 1.599

# ---------------------------------------------------------------------------
# 3) Patch the remaining cells of "Resumen total" that differ from the
#    "Resumen punta" template (everything else is shared between both).
# ---------------------------------------------------------------------------
$totalSheet3 = $wb.Worksheets.Item("Resumen total")
$totalSheet3.Range("D2").Value = 11026.644
$totalSheet3.Range("C4").Value = 13.223
$totalSheet3.Range("C5").Value = 44568.35416666666
$totalSheet3.Range("C6").Value = "intermedia"
$totalSheet3.Range("C7").Value = "08:30:00"

# ---------------------------------------------------------------------------
# 4) Insert the new row 3 into the three pre-existing "Resumen ..." sheets.
# ---------------------------------------------------------------------------
Set-NewRowThree "Resumen Base" 0
Set-NewRowThree "Resumen Intermedia" 1.599
Set-NewRowThree "Resumen punta" 0
